# Generate Report for Handoff
#
# Regenerates the localization-status report for a new handoff: the
# source document got a new GUID-named blob and new target-file commit
# hash, and the handoff timestamps advance a few seconds/minutes.
#
#   ac4567ca-4bb2-4e90-b090-0a0b9cc1e083   ->  92ee2ffd-a547-4a8f-a31c-0bd2f950989f   (source file id)
#   6f912263037c6f4ae72ce2d648f805dd8fdc09d6 -> 465aa127ac377bb7970b35a85db41248b6fbed9e (target commit)

$wb = $excel.ActiveWorkbook

$newGuid = "92ee2ffd-a547-4a8f-a31c-0bd2f950989f"
$newHash = "465aa127ac377bb7970b35a85db41248b6fbed9e"

# ---- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value2 = "$newGuid.md"
$wsOverview.Range("D2").Value2 = "2016-03-22 13:08:14"

# ---- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value2 = "$newGuid.md"
$wsZhCn.Range("D2").Value2 = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value2 = "2016-03-22 13:08:10"

# ---- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value2 = "$newGuid.md"
$wsDeDe.Range("D2").Value2 = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value2 = "2016-03-22 13:08:14"
